$d = $word.ActiveDocument

# Fix 1: "Wendy  de" (double space) -> "Wendy de" (single space)
$d.Content.Find.Execute("Wendy  de", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Wendy de", 2)

# Fix 2: "Jinet" -> "Jinete" (two occurrences)
$d.Content.Find.Execute("Jinet", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Jinete", 2)
